$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Category" header in A1, formatted like the rest of the
# header row (copy B1's format onto A1, then set the text).
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The category-name cells in column A (rows 2-46) previously shared the
# header style; strip that formatting so they fall back to the default
# (unstyled) cell format.
$ws.Range("A2:A46").Style = "Normal"
